$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3, column C (Status) from PASS to FAIL
$ws.Range("C3").Value = "FAIL"

# Add new row 4 duplicating row 3's Test Case ID / Method Name, with FAIL status
$ws.Range("A4").Value = "TC-TESTVALIDLOGIN2"
$ws.Range("B4").Value = "testValidLogin2"
$ws.Range("C4").Value = "FAIL"
